# Natmi following Dr Hou advice
# Re-run of the NATMI ligand(Tgfb2)-receptor(Tgfbr3) pair table now that
# "ECs" is included as a sending cluster as well as a target cluster, so
# the 3 senders (ECs, FAPs, sCs) x 4 targets (ECs, FAPs, M2, sCs) grid
# fully replaces the old 2 senders x 4 targets grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,("ECs", "Tgfb2", "Tgfbr3", "ECs", 2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 20.424575, 61.273725, 0.203732656096709, 0.2037326560967089, 27.06943134039167, 243.624882063525, 0.005945076203068326, 0.005945076203068325)
  ,("ECs", "Tgfb2", "Tgfbr3", "FAPs", 2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 57.16769933333334, 171.503098, 0.5702408607336045, 0.5702408607336045, 75.76642901954244, 681.8978611758821, 0.01664006858849034, 0.01664006858849034)
  ,("ECs", "Tgfb2", "Tgfbr3", "M2", 2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 0.1182046666666667, 0.354614, 0.001179077200040937, 0.001179077200040937, 0.1566609395028889, 1.409948455526, 0.00003440638304060787, 0.00003440638304060787)
  ,("ECs", "Tgfb2", "Tgfbr3", "sCs", 2, 0.6666666666666666, 1.325336333333333, 3.976009, 0.02918077208126263, 0.02918077208126263, 3, 1, 22.54136766666667, 67.62410300000001, 0.2248474059696456, 0.2248474059696456, 29.87489357165856, 268.874042144927, 0.006561220906663359, 0.006561220906663359)
  ,("FAPs", "Tgfb2", "Tgfbr3", "ECs", 3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 20.424575, 61.273725, 0.203732656096709, 0.2037326560967089, 417.7816655960917, 3760.034990364825, 0.09175456281223957, 0.09175456281223955)
  ,("FAPs", "Tgfb2", "Tgfbr3", "FAPs", 3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 57.16769933333334, 171.503098, 0.5702408607336045, 0.5702408607336045, 1169.356848099732, 10524.21163289759, 0.2568179391400585, 0.2568179391400585)
  ,("FAPs", "Tgfb2", "Tgfbr3", "M2", 3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 0.1182046666666667, 0.354614, 0.001179077200040937, 0.001179077200040937, 2.417859001777555, 21.760731015998, 0.0005310180266843499, 0.0005310180266843499)
  ,("FAPs", "Tgfb2", "Tgfbr3", "sCs", 3, 1, 20.45485233333333, 61.364557, 0.4503674794711605, 0.4503674794711605, 3, 1, 22.54136766666667, 67.62410300000001, 0.2248474059696456, 0.2248474059696456, 461.0803470130413, 4149.723123117371, 0.1012639594921781, 0.101263959492178)
  ,("sCs", "Tgfb2", "Tgfbr3", "ECs", 3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 20.424575, 61.273725, 0.203732656096709, 0.2037326560967089, 482.7950690048667, 4345.1556210438, 0.1060330170814011, 0.1060330170814011)
  ,("sCs", "Tgfb2", "Tgfbr3", "FAPs", 3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 57.16769933333334, 171.503098, 0.5702408607336045, 0.5702408607336045, 1351.327180344567, 12161.94462310111, 0.2967828530050556, 0.2967828530050556)
  ,("sCs", "Tgfb2", "Tgfbr3", "M2", 3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 0.1182046666666667, 0.354614, 0.001179077200040937, 0.001179077200040937, 2.794115921630222, 25.147043294672, 0.0006136527903159789, 0.0006136527903159789)
  ,("sCs", "Tgfb2", "Tgfbr3", "sCs", 3, 1, 23.63794933333334, 70.913848, 0.5204517484475769, 0.5204517484475769, 3, 1, 22.54136766666667, 67.62410300000001, 0.2248474059696456, 0.2248474059696456, 532.831706808705, 4795.485361278344, 0.1170222255708042, 0.1170222255708042)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $vals = $data[$i]
  for ($j = 0; $j -lt $vals.Count; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $vals[$j]
  }
}
